# Updated symbol list on Wed Jan 25 18:55:03 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for the coin rows
# whose market data changed in this run. Values are written as literal text
# (leading apostrophe) so Excel doesn't reinterpret numeric-looking strings
# like "0.1900" or "64.67" as numbers and strip the significant trailing
# zeros / percent sign; ClearFormats() afterwards drops the transient
# "quote prefix" cell style Excel applies for text entry, so the cells keep
# their original (unstyled) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "301.04";    "E2"  = "-3.12%"
    "D3"  = "35.47";     "E3"  = "-0.29%"
                         "E4"  = "-0.89%"
    "D5"  = "0.08007";   "E5"  = "-2.43%"
    "D6"  = "1.899";     "E6"  = "-8.14%"
    "D7"  = "7.761";     "E7"  = "-2.45%"
    "D8"  = "0.9268";    "E8"  = "0.22%"
    "D9"  = "0.1450";    "E9"  = "28.12%"
    "D10" = "0.1900";    "E10" = "-0.97%"
    "D11" = "0.08995";   "E11" = "-2.51%"
    "D12" = "0.03430";   "E12" = "-6.77%"
    "D13" = "0.09845";   "E13" = "-0.68%"
    "D14" = "0.001400";  "E14" = "-2.80%"
    "D15" = "0.005822";  "E15" = "-0.18%"
    "D16" = "3.537";     "E16" = "1.64%"
    "D17" = "4.048";     "E17" = "-1.87%"
    "D18" = "2.957";     "E18" = "1.68%"
    "D19" = "0.3444";    "E19" = "1.21%"
                         "E20" = "-0.45%"
    "D21" = "5.044";     "E21" = "-0.98%"
    "D22" = "0.2397";    "E22" = "8.75%"
    "D23" = "0.04484";   "E23" = "-1.25%"
    "D24" = "0.001214";  "E24" = "-0.98%"
    "D25" = "0.004764";  "E25" = "-1.01%"
    "D26" = "0.0001229"; "E26" = "-1.62%"
    "D27" = "0.0003021"; "E27" = "-32.04%"
    "D39" = "0.01835";   "E39" = "-7.36%"
    "D40" = "0.04763";   "E40" = "-2.34%"
    "D41" = "0.01060";   "E41" = "12.39%"
    "D42" = "0.007353";  "E42" = "-3.75%"
    "D43" = "0.1327";    "E43" = "-4.30%"
    "D44" = "0.002108";  "E44" = "-0.58%"
    "D45" = "0.01088";   "E45" = "-6.47%"
    "D46" = "0.00006227";"E46" = "-5.00%"
    "D47" = "0.00000000749"; "E47" = "-0.05%"
    "D48" = "64.67";     "E48" = "-64.11%"
    "D50" = "0.00002098";"E50" = "-0.05%"
    "D51" = "0.0001998"; "E51" = "-0.05%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.ClearFormats()
}
